$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.589
$ws.Range("A6").Value = -22.225
$ws.Range("A7").Value = -19.575
$ws.Range("C7").Value = -12.489
$ws.Range("C12").Value = -10.782
$ws.Range("C15").Value = -13.447
$ws.Range("A16").Value = -21.857
$ws.Range("A20").Value = -19.853
$ws.Range("C20").Value = -12.366
$ws.Range("C21").Value = -12.519
$ws.Range("C22").Value = -12.867
$ws.Range("C23").Value = -12.223
$ws.Range("A28").Value = -21.942
$ws.Range("A29").Value = -21.426
$ws.Range("C29").Value = -12.305
$ws.Range("A32").Value = -21.771
$ws.Range("C34").Value = -11.955
$ws.Range("A40").Value = -20.067
$ws.Range("C42").Value = -12.492
$ws.Range("C43").Value = -13.159
$ws.Range("C44").Value = -13.216
$ws.Range("C45").Value = -13.054
$ws.Range("A46").Value = -21.951
$ws.Range("C46").Value = -13.862
$ws.Range("C50").Value = -14.19
$ws.Range("A51").Value = -21.697
$ws.Range("C51").Value = -11.2
$ws.Range("A52").Value = -21.907
$ws.Range("A57").Value = -22.247
$ws.Range("A59").Value = -22.427
$ws.Range("A62").Value = -22.158
$ws.Range("A66").Value = -21.674
$ws.Range("C66").Value = -11.405
$ws.Range("C67").Value = -11.3
$ws.Range("A73").Value = -20.597
$ws.Range("A74").Value = -21.244
$ws.Range("C79").Value = -11.816
$ws.Range("C84").Value = -14.098
$ws.Range("A92").Value = -21.626
$ws.Range("C92").Value = -11.472
$ws.Range("C97").Value = -12.788
$ws.Range("A100").Value = -22.063
